$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Nikk Dwivedi"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "8319693453"
$ws.Range("D2").Value = "I'm a Developer"
$ws.Range("E2").Value = "angular"
